{"js": "// Highlight the two \"point system\" bullets green and add a new bullet\n// (\"COMMENT THE CODE THOUROUGLY!\") right after them, per the commit:\n// \"Added a wee bit of a point system\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet tempLocationPara = null;\nlet shoddyMenuPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Add points to a temporary location\") !== -1) {\n    tempLocationPara = p;\n  } else if (p.text.indexOf(\"This should be done in a very shoddy menu\") !== -1) {\n    shoddyMenuPara = p;\n  }\n}\nif (!tempLocationPara || !shoddyMenuPara) {\n  throw new Error(\"Could not locate the target paragraphs.\");\n}\n\n// Insert the new list bullet first (it inherits ilvl=2 from its anchor\n// paragraph) so the highlight applied below doesn't bleed onto it, then\n// reset it back to the top list level (ilvl 0) to match its siblings.\nconst newPara = shoddyMenuPara.insertParagraph(\"COMMENT THE CODE THOUROUGLY!\", \"After\");\nnewPara.listItem.level = 0;\n\n// Green-highlight both the paragraph mark and run text of the two bullets.\ntempLocationPara.font.highlightColor = \"BrightGreen\";\nshoddyMenuPara.font.highlightColor = \"BrightGreen\";\n\nawait context.sync();\n", "ps1": "# Highlight the two \"point system\" bullets green and add a new bullet\n# (\"COMMENT THE CODE THOUROUGLY!\") right after them, per the commit:\n# \"Added a wee bit of a point system\".\n$d = $word.ActiveDocument\n\n$target1 = $null\n$target2 = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Add points to a temporary location*\") {\n        $target1 = $p\n    } elseif ($p.Range.Text -like \"*This should be done in a very shoddy menu*\") {\n        $target2 = $p\n    }\n}\nif ($target1 -eq $null -or $target2 -eq $null) {\n    throw \"Could not locate the target paragraphs.\"\n}\n\n# Insert the new list bullet first (it inherits ilvl=2/numId=3 from its\n# anchor paragraph) so the highlight applied below doesn't bleed onto it,\n# then reset it back to the top list level (wdListLevelNumber 1 == ilvl 0)\n# to match its siblings.\n$null = $target2.Range.InsertParagraphAfter()\n$newPara = $target2.Next()\n$newPara.Range.Text = \"COMMENT THE CODE THOUROUGLY!\"\n$newPara.Range.ListFormat.ListLevelNumber = 1\n\n# Green-highlight both the paragraph mark and run text of the two bullets.\n$target1.Range.Font.HighlightColorIndex = \"wdBrightGreen\"\n$target2.Range.Font.HighlightColorIndex = \"wdBrightGreen\"\n"}
